{"js": "// The document is a single table, one column, where each row holds a\n// single benchmark value in its one cell. This edit refreshes the\n// \"summary\" rows (1-12, 0-indexed 0-11) with corrected/updated stats and\n// collapses the verbose per-run tab-separated rows (rows 44-46, 0-indexed\n// 43-45) down to the single overall value that used to live in rows 1-3.\n\nconst table = context.document.body.tables.getFirst();\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Need the cells collection loaded for every row we touch before we can\n// reach `.items[0]` on each of them.\nfor (let i = 0; i < rows.items.length; i++) {\n  rows.items[i].cells.load(\"items\");\n}\nawait context.sync();\n\nfunction setRowValue(rowIndex, newValue) {\n  rows.items[rowIndex].cells.items[0].value = newValue;\n}\n\n// Simple value replacements (single run per cell already).\nsetRowValue(0, \"0M\");\nsetRowValue(1, \"0M\");\nsetRowValue(2, \"0M\");\nsetRowValue(3, \"1217\");\nsetRowValue(5, \"0.00077\");\nsetRowValue(6, \"0.00020\");\nsetRowValue(7, \"0.00006\");\nsetRowValue(8, \"0.00029\");\nsetRowValue(9, \"0.00035\");\nsetRowValue(10, \"0.00045\");\nsetRowValue(11, \"0.24014\");\n\n// Collapse the tab-separated detail rows into the single summary value\n// they used to hold (this also drops the extra <w:t>/<w:tab/> runs).\nsetRowValue(43, \"99.73\");\nsetRowValue(44, \"0.24\");\nsetRowValue(45, \"88\");\n\nawait context.sync();\n", "ps1": "# The document is a single table, one column, where each row holds a\n# single benchmark value in its one cell. This edit refreshes the\n# \"summary\" rows (1-12) with corrected/updated stats and collapses the\n# verbose per-run tab-separated rows (44-46) down to the single overall\n# value that used to live in rows 1-3.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Simple value replacements (single run per cell already).\n$t.Cell(1, 1).Range.Text = \"0M\"\n$t.Cell(2, 1).Range.Text = \"0M\"\n$t.Cell(3, 1).Range.Text = \"0M\"\n$t.Cell(4, 1).Range.Text = \"1217\"\n$t.Cell(6, 1).Range.Text = \"0.00077\"\n$t.Cell(7, 1).Range.Text = \"0.00020\"\n$t.Cell(8, 1).Range.Text = \"0.00006\"\n$t.Cell(9, 1).Range.Text = \"0.00029\"\n$t.Cell(10, 1).Range.Text = \"0.00035\"\n$t.Cell(11, 1).Range.Text = \"0.00045\"\n$t.Cell(12, 1).Range.Text = \"0.24014\"\n\n# Collapse the tab-separated detail rows into the single summary value\n# they used to hold (this also drops the extra text runs/tab chars).\n$t.Cell(44, 1).Range.Text = \"99.73\"\n$t.Cell(45, 1).Range.Text = \"0.24\"\n$t.Cell(46, 1).Range.Text = \"88\"\n"}
